$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "enabled" column (D, rows 2-7) used to store native Boolean TRUE
# values. It is now stored as the literal text "TRUE" (a shared string)
# instead, with an in-cell dropdown data validation restricting entry to
# TRUE/FALSE.
#
# A plain Range.Value = "TRUE" assignment gets re-interpreted as a native
# Boolean by the engine's type inference (same as typing TRUE into Excel),
# so instead we compute the text via TEXT(TRUE,""), then paste-special
# "values only" over itself so the result lands as a plain string cell
# (no quote-prefix / no extra cell style, matching a programmatically
# generated import file).
$boolRange = $ws.Range("D2:D7")
$boolRange.Formula = '=TEXT(TRUE,"")'
$boolRange.Copy()
$boolRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Add a list data validation on the whole enabled column (minus header),
# matching Excel's default "apply to whole column" behavior.
$validatedRange = $ws.Range("D2:D1048576")
$validatedRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$validatedRange.Validation.ErrorTitle = "Enabled Error"
$validatedRange.Validation.ErrorMessage = "You must choose true or false"
$validatedRange.Validation.IgnoreBlank = $true
$validatedRange.Validation.InCellDropdown = $true
$validatedRange.Validation.ShowInput = $true
$validatedRange.Validation.ShowError = $true
